$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.770179333333333
$ws.Range("H2").Value = 5.310538
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("Q2").Value = 21.87246447569822
$ws.Range("R2").Value = 196.852180281284
